# Apply NBA MVP pivot table data cleanup:
#  - Insert a new row for "Joel Embiid" (1 award) in alphabetical order
#    between "James Harden" and "Julius Erving".
#  - Correct "LeBron James" award count (was incorrectly 10, should be 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "James Harden" currently sits on row 15, "Julius Erving" on row 16.
# Insert a fresh row above row 16 (Julius Erving) to hold "Joel Embiid",
# shifting Julius Erving and everyone after it down by one row.
$insertRow = 16
$ws.Rows.Item($insertRow).Insert()

$ws.Cells.Item($insertRow, 1).Value = "Joel Embiid"
$ws.Cells.Item($insertRow, 2).Value = 1

# Fix LeBron James' award total, which is now a few rows further down
# because of the inserted row above.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $playerName = $ws.Cells.Item($r, 1).Value()
    if ($playerName -eq "LeBron James") {
        $ws.Cells.Item($r, 2).Value = 4
        break
    }
}
